# A new weekly price observation is inserted into the "Ciboulette" price
# series for "Vega Modelo de Temuco". The new record belongs right above the
# current row 147 (it shares that row's price/volume-unit/origin data, just
# with a later date and a different reported volume), so every existing
# record from row 147 down shifts one row further down the sheet
# (147-169 -> 148-170) and the used range grows from A1:R169 to A1:R170.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 147 (and everything below it) down by one row, opening up a
# blank row 147 for the new record.
$ws.Rows(147).Insert()

# Seed the new row 147 with the same data as the record that is now sitting
# in row 148 (the former row 147) so every column besides the date/volume
# starts out identical, then overwrite the two columns that actually differ
# for the new observation (Fecha / Volumen).
$copied = $ws.Range("A148:R148").Value()
$ws.Range("A147:R147").Value = $copied

$ws.Cells.Item(147, 4).Value = 44504   # Fecha -> 2021-11-04
$ws.Cells.Item(147, 10).Value = 125    # Volumen -> 125
